$d = $word.ActiveDocument
$x = [char]0xD7

# Update the date line
$d.Content.Find.Execute("2026-02-13 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-14 Saturday", 2)

$t = $d.Tables.Item(1)

# Mapping of (row, col) -> (old, new) values, using the multiplication sign $x
$edits = @(
    @{ Row = 1;  Col = 1; Old = "295${x}2="; New = "757${x}5=" },
    @{ Row = 1;  Col = 2; Old = "113${x}6="; New = "223${x}4=" },
    @{ Row = 1;  Col = 3; Old = "834${x}3="; New = "329${x}9=" },
    @{ Row = 1;  Col = 4; Old = "647${x}3="; New = "932${x}5=" },
    @{ Row = 1;  Col = 5; Old = "157${x}6="; New = "360${x}8=" },

    @{ Row = 5;  Col = 1; Old = "634${x}7="; New = "937${x}4=" },
    @{ Row = 5;  Col = 2; Old = "841${x}7="; New = "762${x}3=" },
    @{ Row = 5;  Col = 3; Old = "814${x}4="; New = "923${x}8=" },
    @{ Row = 5;  Col = 4; Old = "209${x}4="; New = "582${x}6=" },
    @{ Row = 5;  Col = 5; Old = "481${x}7="; New = "959${x}3=" },

    @{ Row = 10; Col = 1; Old = "439${x}5="; New = "383${x}2=" },
    @{ Row = 10; Col = 2; Old = "359${x}9="; New = "518${x}8=" },
    @{ Row = 10; Col = 3; Old = "226${x}9="; New = "106${x}9=" },
    @{ Row = 10; Col = 4; Old = "983${x}4="; New = "158${x}4=" },
    @{ Row = 10; Col = 5; Old = "586${x}7="; New = "258${x}2=" },

    @{ Row = 15; Col = 1; Old = "683${x}2="; New = "569${x}2=" },
    @{ Row = 15; Col = 2; Old = "583${x}8="; New = "673${x}7=" },
    @{ Row = 15; Col = 3; Old = "440${x}6="; New = "330${x}8=" },
    @{ Row = 15; Col = 4; Old = "897${x}7="; New = "172${x}7=" },
    @{ Row = 15; Col = 5; Old = "683${x}2="; New = "818${x}4=" },

    @{ Row = 20; Col = 1; Old = "516${x}8="; New = "986${x}6=" },
    @{ Row = 20; Col = 2; Old = "677${x}6="; New = "453${x}6=" },
    @{ Row = 20; Col = 3; Old = "649${x}3="; New = "969${x}5=" },
    @{ Row = 20; Col = 4; Old = "798${x}3="; New = "440${x}8=" },
    @{ Row = 20; Col = 5; Old = "913${x}3="; New = "149${x}6=" }
)

foreach ($edit in $edits) {
    $cellRange = $t.Cell($edit.Row, $edit.Col).Range
    $rng = $d.Range($cellRange.Start, $cellRange.End)
    $ok = $rng.Find.Execute($edit.Old, $true, $false, $false, $false, $false, $true, 1, $false, $edit.New, 2)
    if (-not $ok) {
        Write-Host "FAILED:" $edit.Row $edit.Col $edit.Old "->" $edit.New
    }
}

Write-Host "Done"
